$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# 1.1 intro paragraph: append "的地方" right before the closing period.
$ok1 = $d.Content.Find.Execute(
    "因此，将可信计算技术应用在提高云计算环境的安全性是工业界和产业界必须重视。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "因此，将可信计算技术应用在提高云计算环境的安全性是工业界和产业界必须重视的地方。",
    2)
if (-not $ok1) { throw "Change 1: target sentence not found" }

# --- Change 2 -------------------------------------------------------------
# The paragraph with the second <w:commentReference> gets a big rewrite:
#   a) the tail of the existing (pre-bookmark) run is replaced with new
#      wording that now ends right before the "_GoBack" bookmark, and
#   b) a brand-new run is appended immediately *after* the bookmark with
#      its own sentence, using the same run formatting
#      (ascii=Times New Roman / hAnsi=宋体 / hint=eastAsia / szCs=21 /
#      lang en-US,eastAsia=zh-CN) as the run that precedes it.

# 2a. Replace the old closing sentence (still inside the existing run,
#     before the bookmark) with the new closing wording.
$ok2 = $d.Content.Find.Execute(
    "而可信虚拟平台的构建可以利用TPM中的可信度量、可信报告等技术向用户发送关于云计算平台的可信度量结果，并且证明自身的安全性。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "主要是由于云租户在使用云提供商提供的虚拟机时，并不能确认云计算平台上的物理主机是与云提供商按照各自操作系统官方文件进行启动的，以及租户请求的虚拟机是按照预期的配置和要求进行启动的。因为云计算环境下的虚拟机存在着包括传统信息系统安全以及新型网络安全等威胁，比如：虚拟机恶意代码攻击、虚拟机逃逸等，这些都会导致虚拟机在重新启动时的组件被篡改，",
    2)
if (-not $ok2) { throw "Change 2a: target sentence not found" }

# 2b. Insert the new trailing run right after the "_GoBack" bookmark.
#     Trick: temporarily stretch a tiny range that already carries the
#     correct run formatting (the character immediately in front of the
#     bookmark) to hold the new sentence, copy it (with its formatting)
#     to the insertion point after the bookmark via FormattedText, then
#     restore the original range back to its original content.
$newSentence = "在云租户对虚拟机进行重新启动时，可能无法判断虚拟机遭受操作系统、数据是否被篡改。而可信虚拟平台的构建可以利用TPM中的可信度量、可信报告等技术向用户发送关于云计算平台的可信度量结果，并且证明自身的安全性。"

$bm = $d.Bookmarks.Item("_GoBack")
$srcRange = $d.Range($bm.Start - 1, $bm.Start)
$origText = $srcRange.Text

$srcRange.Text = $newSentence
$destRange = $d.Range($srcRange.End, $srcRange.End)
$destRange.FormattedText = $srcRange.FormattedText

$srcRange.Text = $origText
